$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 111651416
$ws.Range("I2").Value = "'10"
$ws.Range("I2").Style = "Normal"
$ws.Range("Q2").Value = 573987
$ws.Range("R2").Value = 6403999
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# --- Row 3 ---
$ws.Range("A3").Value = 111651333
$ws.Range("I3").Value = "'25"
$ws.Range("I3").Style = "Normal"
$ws.Range("Q3").Value = 574013
$ws.Range("R3").Value = 6403975
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# --- Row 4 ---
$ws.Range("Q4").Value = 574026
$ws.Range("R4").Value = 6403965
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()

# --- Row 5 ---
$ws.Range("Q5").Value = 574016
$ws.Range("R5").Value = 6403975

# --- Row 6 ---
$ws.Range("Q6").Value = 574025
$ws.Range("R6").Value = 6403972
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
